# Applies the "Updated cryptos list on Mon Aug 28 16:52:41 UTC 2023 with GitHub Actions"
# edit: refreshed Price/Volume(1h) figures and swapped the Frax/BabyDogeCoin rows (48-49).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells (column D) that parse as plain numbers need to be forced back to
# Text format first, otherwise Excel auto-converts the typed string into a float
# (e.g. "0.5237" -> 0.52370000000000005) and the cell loses its inlineStr/text type.
$numericLookingPriceCells = @(
    "D5",
    "D6",
    "D8",
    "D9",
    "D10",
    "D11",
    "D12",
    "D15",
    "D17",
    "D20",
    "D21",
    "D22",
    "D23",
    "D24",
    "D25",
    "D26",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D34",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D45",
    "D47",
    "D48",
    "D50"
)
foreach ($cellRef in $numericLookingPriceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.239.51"
$ws.Range("E2").Value = "  -0.69%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.656.98"
$ws.Range("E3").Value = "  -0.97%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.59%  "

# Row 5 - BNB
$ws.Range("D5").Value = "219.84"
$ws.Range("E5").Value = "  -0.85%  "

# Row 6 - XRP
$ws.Range("D6").Value = "0.5237"
$ws.Range("E6").Value = "  -1.83%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.57%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "0.2675"
$ws.Range("E8").Value = "  +0.41%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "0.06370"
$ws.Range("E9").Value = "  -0.36%  "

# Row 10 - Solana
$ws.Range("D10").Value = "20.65"
$ws.Range("E10").Value = "  -1.39%  "

# Row 11 - TRON
$ws.Range("D11").Value = "0.07722"
$ws.Range("E11").Value = "  -1.55%  "

# Row 12 - Polkadot
$ws.Range("D12").Value = "4.613"
$ws.Range("E12").Value = "  +1.61%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.633.91"
$ws.Range("E13").Value = "  -2.92%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "1.884.79"
$ws.Range("E14").Value = "  -0.85%  "

# Row 15 - Polygon
$ws.Range("D15").Value = "0.5658"
$ws.Range("E15").Value = "  +0.55%  "

# Row 16 - ShibaInu
$ws.Range("D16").Value = "0.0₅8253"
$ws.Range("E16").Value = "  +0.78%  "

# Row 17 - Litecoin
$ws.Range("D17").Value = "65.47"
$ws.Range("E17").Value = "  -1.16%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "26.236.87"
$ws.Range("E18").Value = "  -0.63%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  -0.62%  "

# Row 20 - Uniswap
$ws.Range("D20").Value = "4.708"
$ws.Range("E20").Value = "  -0.31%  "

# Row 21 - Avalanche
$ws.Range("D21").Value = "10.42"
$ws.Range("E21").Value = "  +0.99%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "192.26"
$ws.Range("E22").Value = "  -2.89%  "

# Row 23 - Chainlink
$ws.Range("D23").Value = "6.018"
$ws.Range("E23").Value = "  -1.10%  "

# Row 24 - BinanceUSD
$ws.Range("D24").Value = "1.005"
$ws.Range("E24").Value = "  -0.51%  "

# Row 25 - Monero
$ws.Range("D25").Value = "143.24"
$ws.Range("E25").Value = "  -2.36%  "

# Row 26 - Stellar
$ws.Range("D26").Value = "0.1202"
$ws.Range("E26").Value = "  -2.23%  "

# Row 27 - Cosmos
$ws.Range("D27").Value = "7.290"
$ws.Range("E27").Value = "  +0.56%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "15.98"
$ws.Range("E28").Value = "  -1.86%  "

# Row 29 - Toncoin
$ws.Range("D29").Value = "1.510"
$ws.Range("E29").Value = "  +0.32%  "

# Row 30 - Hedera
$ws.Range("D30").Value = "0.05645"
$ws.Range("E30").Value = "  -4.59%  "

# Row 31 - PancakeSwap
$ws.Range("D31").Value = "1.278"
$ws.Range("E31").Value = "  -0.95%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("D32").Value = "3.509"
$ws.Range("E32").Value = "  -1.53%  "

# Row 33 - Filecoin
$ws.Range("E33").Value = "  +1.08%  "

# Row 34 - LidoDAOToken
$ws.Range("D34").Value = "1.580"
$ws.Range("E34").Value = "  -2.34%  "

# Row 35 - MXToken
$ws.Range("E35").Value = "  -1.26%  "

# Row 36 - ARBITRUM
$ws.Range("D36").Value = "0.9461"
$ws.Range("E36").Value = "  -2.57%  "

# Row 37 - HuobiToken
$ws.Range("D37").Value = "2.412"
$ws.Range("E37").Value = "  -0.93%  "

# Row 38 - ImmutableX
$ws.Range("D38").Value = "0.5781"
$ws.Range("E38").Value = "  -0.95%  "

# Row 39 - VeChain
$ws.Range("D39").Value = "0.01602"
$ws.Range("E39").Value = "  -0.85%  "

# Row 40 - FraxShare
$ws.Range("D40").Value = "5.925"
$ws.Range("E40").Value = "  +0.34%  "

# Row 41 - mCoin
$ws.Range("D41").Value = "2.569"
$ws.Range("E41").Value = "  -0.37%  "

# Row 42 - TrustWalletToken
$ws.Range("D42").Value = "0.8478"
$ws.Range("E42").Value = "  -1.96%  "

# Row 44 - Maker
$ws.Range("D44").Value = "1.023.04"
$ws.Range("E44").Value = "  -5.36%  "

# Row 45 - Quant
$ws.Range("D45").Value = "101.48"
$ws.Range("E45").Value = "  -2.21%  "

# Row 46 - RocketPoolETH (price only; Volume(1h) unchanged)
$ws.Range("D46").Value = "1.795.03"

# Row 47 - Aave
$ws.Range("D47").Value = "58.50"
$ws.Range("E47").Value = "  -0.06%  "

# Row 48 - was BabyDogeCoin, now Frax (rows 48/49 swapped places)
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").Value = "1.006"
$ws.Range("E48").Value = "  -0.47%  "

# Row 49 - was Frax, now BabyDogeCoin
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₈104"
$ws.Range("E49").Value = "  -1.90%  "

# Row 50 - Cronos
$ws.Range("D50").Value = "0.05324"
$ws.Range("E50").Value = "  +3.07%  "

# Row 51 - Mantle
$ws.Range("E51").Value = "  -1.37%  "

